# Feria Lagunitas de Puerto Montt - Pepino ensalada
# Adds a new weekly price record, inserted as row 166 (pushing the
# existing historical rows 166-215 down to 167-216).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 166; existing rows 166-215 move to 167-216.
$ws.Rows("166:166").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(166, 1).Value = 4
$ws.Cells.Item(166, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(166, 3).Value = "Los Lagos"
$ws.Cells.Item(166, 4).Value = 44588
$ws.Cells.Item(166, 5).Value = 10
$ws.Cells.Item(166, 6).Value = 100112043
$ws.Cells.Item(166, 7).Value = "Pepino ensalada"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 200
$ws.Cells.Item(166, 11).Value = 13000
$ws.Cells.Item(166, 12).Value = 14000
$ws.Cells.Item(166, 13).Value = 13500
$ws.Cells.Item(166, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(166, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(166, 16).Value = 225
$ws.Cells.Item(166, 17).Value = 60
$ws.Cells.Item(166, 18).Value = "Hortaliza"
